$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 32: "Wachen-Sichtradius schlägt jetzt auch schon an..." entry.
# Write the long description (column E) before the version/testmethod/tester
# cells so the shared-string table append order matches the target file
# (new description string lands before the new "...exe" version string).
$ws.Range("E32").Value = "Wachen-Sichtradius schlägt jetzt auch schon an, wenn man seitlich rein kommt, d.h. wenn noch nicht Mittelpunkt des Spielers oder toten Körpers drin ist"
$ws.Range("B32").Value = "DiscordiaAgency_Demo_2017_09_23-3.exe"
$ws.Range("C32").Value = "Entwicklung"
$ws.Range("D32").Value = "Anna Franziska"

# Match the row height used for the other wrapped-text entries in this column.
$ws.Rows.Item(32).RowHeight = 45

# Move the selection the way it ended up after the edit in the source workbook.
$ws.Range("E38").Select()
